$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 3068
$ws1.Range("F6").Value = 2050
$ws1.Range("F8").Value = 144
$ws1.Range("F9").Value = 1148
$ws1.Range("F10").Value = 210
$ws1.Range("F11").Value = 871
$ws1.Range("F12").Value = 73

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 3068
$ws4.Range("F6").Value = 2050
$ws4.Range("F9").Value = 144
$ws4.Range("F10").Value = 1148
$ws4.Range("F11").Value = 210
$ws4.Range("F12").Value = 871
$ws4.Range("F13").Value = 73
